$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need to be forced to
# text (format as Text, assign, then clear the format again) so Excel
# does not silently convert them into numeric cells - the source data
# keeps these as strings (e.g. "564.53"), matching the original file.

$ws.Range("D2").Value = "60.656.24"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "2.400.64"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.53"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  +1.46%  "
$ws.Range("D9").Value = "2.407.55"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.19"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.06"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("D16").Value = "2.810.98"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "60.344.69"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").Value = "2.416.60"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.08"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +8.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.64"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.13"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.04"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.83"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "568.45"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.04"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.14%  "
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.08"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.23%  "
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.132"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").ClearFormats()
$ws.Range("E36").Value = "  +4.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.23"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.45%  "
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.31"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.73"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.67"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("E45").Value = "  +6.16%  "
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.06"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0507"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.31"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.06%  "
